# Daily attendance processing - 2026-02-14 15:32:36 UTC
# Swap the order of "Administrator, Miss Dina Nasr" -> "Miss Dina Nasr, Administrator"
# in the "Recorded By" column (column G) of the attendance sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldText = "Administrator, Miss Dina Nasr"
$newText = "Miss Dina Nasr, Administrator"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    if ($cell.Value2 -eq $oldText) {
        $cell.Value = $newText
    }
}
